$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 13
$ws.Range("B4").Value = 45
$ws.Range("B5").Value = 46
$ws.Range("B6").Value = 100
$ws.Range("B7").Value = 64
$ws.Range("B8").Value = 83
$ws.Range("B9").Value = 245
$ws.Range("B10").Value = 543
$ws.Range("B11").Value = 729
$ws.Range("B12").Value = 828
$ws.Range("B13").Value = 1616
$ws.Range("B14").Value = 2639
$ws.Range("B15").Value = 2340
$ws.Range("B16").Value = 2355
